# Refactor workbook: drop the "customer order" sheet (duplicate of
# "supplier order") and repurpose the remaining one as "customer deadlines".

$wb = $excel.ActiveWorkbook

# Suppress the "delete sheet" confirmation alert, if available.
if ($excel.PSObject.Properties.Match('DisplayAlerts').Count -gt 0) {
    $excel.DisplayAlerts = $false
}

# Remove the now-redundant "customer order" sheet entirely.
$custOrder = $wb.Worksheets.Item("customer order")
$custOrder.Delete()

# Rename "supplier order" -> "customer deadlines" (data/layout unchanged).
$supplierOrder = $wb.Worksheets.Item("supplier order")
$supplierOrder.Name = "customer deadlines"
